$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2:A5 with the new combined tuple-style strings
$ws.Range("A2").Value = "('Manifest', ['Creature', '(You can cover a face-down manifested creature with this reminder card.', 'A manifested creature card can be turned face up any time for its mana cost. A face-down card can also be turned face up for its morph cost.)', '2/2'])"
$ws.Range("A3").Value = "('Monk', ['Token Creature — Monk', 'Prowess', '1/1'])"
$ws.Range("A4").Value = "('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])"
$ws.Range("A5").Value = "('Warrior', ['Token Creature — Warrior', '2/1'])"

# Remove the now-obsolete rows 6 through 17
$ws.Range("A6:A17").ClearContents() | Out-Null
